$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Some price strings (e.g. "210.50", "7.09") look like valid numeric
    # literals, and a plain .Value assignment would have Excel silently
    # convert them to numbers (losing formatting like trailing zeros).
    # Forcing the cell to Text before the write keeps it a real string;
    # ClearFormats() afterwards drops the temporary number-format override
    # again so the cell is left with no explicit style, matching the rest
    # of the sheet.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).ClearFormats()
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.630.30"
$ws.Range("E2").Value = "  -0.30%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.596.47"
$ws.Range("E3").Value = "  -0.24%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.08%  "

# Row 5 - BNB
Set-TextValue "D5" "210.50"
$ws.Range("E5").Value = "  -0.45%  "

# Row 6 - XRP
Set-TextValue "D6" "0.509"
$ws.Range("E6").Value = "  -0.61%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.10%  "

# Row 8 - Dogecoin
$ws.Range("E8").Value = "  -0.69%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -0.55%  "

# Row 10 - Solana
Set-TextValue "D10" "19.54"
$ws.Range("E10").Value = "  +0.01%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.12%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.820.18"
$ws.Range("E12").Value = "  -0.25%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.592.69"
$ws.Range("E13").Value = "  -0.50%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -0.16%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  -0.28%  "

# Row 16 - Litecoin
Set-TextValue "D16" "64.54"
$ws.Range("E16").Value = "  -1.26%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "26.597.95"
$ws.Range("E17").Value = "  -0.33%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "0.0₃0738"

# Row 19 - Dai
$ws.Range("E19").Value = "  +0.09%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "208.00"
$ws.Range("E20").Value = "  -0.71%  "

# Row 21 - Chainlink
Set-TextValue "D21" "7.09"
$ws.Range("E21").Value = "  -1.44%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -0.28%  "

# Row 23 - Toncoin
$ws.Range("E23").Value = "  -3.33%  "

# Row 25 - Monero
Set-TextValue "D25" "143.76"
$ws.Range("E25").Value = "  +0.62%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  +0.05%  "

# Row 27 - Cosmos
$ws.Range("E27").Value = "  +0.05%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  -0.97%  "

# Row 29 - EthereumClassic
$ws.Range("E29").Value = "  -0.58%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  -2.27%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -0.43%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -0.28%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextValue "D33" "2.95"
$ws.Range("E33").Value = "  -0.43%  "

# Row 34/35 swap: WEMIXToken <-> Maker (with updated data)
$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D34").Value = "1.279.49"
$ws.Range("E34").Value = "  -1.01%  "

$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D35" "1.25"
$ws.Range("E35").Value = "  +16.04%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  +0.53%  "

# Row 37 - LidoDAOToken
$ws.Range("E37").Value = "  -1.15%  "

# Row 38 - ImmutableX
Set-TextValue "D38" "0.595"
$ws.Range("E38").Value = "  -4.01%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  -2.19%  "

# Row 40 - ARBITRUM
$ws.Range("E40").Value = "  -0.52%  "

# Row 41 - FraxShare
$ws.Range("E41").Value = "  +0.27%  "

# Row 42 - MXToken
Set-TextValue "D42" "2.15"
$ws.Range("E42").Value = "  -1.50%  "

# Row 43 - TrustWalletToken
Set-TextValue "D43" "0.771"

# Row 44 - Aave
Set-TextValue "D44" "62.53"
$ws.Range("E44").Value = "  -0.98%  "

# Row 45 - RocketPoolETH
$ws.Range("D45").Value = "1.731.74"
$ws.Range("E45").Value = "  -0.20%  "

# Row 46 - Quant
Set-TextValue "D46" "89.41"
$ws.Range("E46").Value = "  -1.66%  "

# Row 48 - Algorand
Set-TextValue "D48" "0.102"
$ws.Range("E48").Value = "  +1.84%  "

# Row 49 - Cronos
$ws.Range("E49").Value = "  +0.60%  "

# Row 50 - USDD
$ws.Range("E50").Value = "  +0.09%  "

# Row 51 - EnergySwap
Set-TextValue "D51" "7.43"
$ws.Range("E51").Value = "  +0.93%  "
